$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 156, pushing existing rows 156-212 down to 157-213.
$ws.Rows("156:156").Insert()

# Fill in the new row 156 with the "Political Regimes of the World" record.
$ws.Range("A156").Value = "Political Regimes of the World"
$ws.Range("B156").Value = "democracy"
$ws.Range("C156").Value = "https://link.springer.com/article/10.1057/s41304-018-0149-8"
$ws.Range("D156").Value = "autocracy, democracy, political regimes"
$ws.Range("E156").Value = "world"
$ws.Range("F156").Value = 1
$ws.Range("G156").Value = 1
$ws.Range("H156").Value = 1
$ws.Range("I156").Value = 1
$ws.Range("J156").Value = 1
$ws.Range("K156").Value = 1800
$ws.Range("L156").Value = 2016
$ws.Range("M156").Value = "online"
$ws.Range("N156").Value = "no"
$ws.Range("O156").Value = 1
$ws.Range("P156").Value = "https://static-content.springer.com/esm/art%3A10.1057%2Fs41304-018-0149-8/MediaObjects/41304_2018_149_MOESM1_ESM.docx"
$ws.Range("T156").Value = "https://static-content.springer.com/esm/art%3A10.1057%2Fs41304-018-0149-8/MediaObjects/41304_2018_149_MOESM2_ESM.xlsx"
$ws.Range("W156").Value = "country"
$ws.Range("X156").Value = "year"
$ws.Range("Y156").Value = "ccode"
$ws.Range("Z156").Value = "10.1057/s41304-018-0149-8"
$ws.Range("AA156").Value = "10.1057/s41304-018-0149-8"
$ws.Range("AB156").Value = 20180320

# Wire up the hyperlinks for the link / file_codebook / file_excel cells.
$ws.Hyperlinks.Add($ws.Range("C156"), "https://link.springer.com/article/10.1057/s41304-018-0149-8")
$ws.Hyperlinks.Add($ws.Range("P156"), "https://static-content.springer.com/esm/art%3A10.1057%2Fs41304-018-0149-8/MediaObjects/41304_2018_149_MOESM1_ESM.docx")
$ws.Hyperlinks.Add($ws.Range("T156"), "https://static-content.springer.com/esm/art%3A10.1057%2Fs41304-018-0149-8/MediaObjects/41304_2018_149_MOESM2_ESM.xlsx")

# Restore the plain "Hyperlink" cell style (Excel's Hyperlinks.Add can otherwise
# spawn a duplicate style) and make sure unrelated cells in the row keep no style.
$ws.Range("C156").Style = "Hyperlink"
$ws.Range("P156").Style = "Hyperlink"
$ws.Range("T156").Style = "Hyperlink"

# Row-insert copies formatting down from row 155, which leaves a stray styled
# (but empty) V156 cell behind (row 155 has a hyperlink in V). Drop it.
$ws.Range("V156").Clear()

# Match the view state captured in the saved workbook (scrolled to the new row).
$ws.Range("A156").Select()
